$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "Yes" -> "No"
$ws.Range("B2").Value = "No"

# F2: "Na" -> empty (clear contents)
$ws.Range("F2").ClearContents()

# J2: 1 -> 0
$ws.Range("J2").Value = 0
